$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Delete all data rows, keeping only the header row (row 1).
for ($i = $table.Rows.Count; $i -ge 2; $i--) {
    $table.Rows.Item($i).Delete()
}
